# Update data: 2025-10-29 18:21

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump "Last Updated" timestamp ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "29 Oct 2025, 06:21 PM"

# --- Top Losers sheet ---
$wsLosers = $wb.Worksheets.Item("Top Losers")

# Row 18: Weekly value update (CRAMC)
$wsLosers.Range("D18").Value = 5.978

# Rows 35-36: SPARC / PRUDENT reorder with refreshed figures
$wsLosers.Range("B35").Value = "SPARC"
$wsLosers.Range("C35").Value = -3.1709
$wsLosers.Range("D35").Value = 4.8337
$wsLosers.Range("E35").Value = 6.3311

$wsLosers.Range("B36").Value = "PRUDENT"
$wsLosers.Range("C36").Value = -3.127
$wsLosers.Range("D36").Value = -3.5103
$wsLosers.Range("E36").Value = 2.1213

# Row 48: Weekly value update (RUBICON)
$wsLosers.Range("D48").Value = -2.9654

# Row 54: Weekly value becomes N/A (CANHLIFE)
$wsLosers.Range("D54").Value = "N/A"

# Rows 64-76: refreshed list (values shift up one row, new entry BBOX added at bottom)
$wsLosers.Range("B64").Value = "NESCO"
$wsLosers.Range("C64").Value = -2.4722
$wsLosers.Range("D64").Value = 1.9934
$wsLosers.Range("E64").Value = 3.8931

$wsLosers.Range("B65").Value = "PILANIINVS"
$wsLosers.Range("C65").Value = -2.4546
$wsLosers.Range("D65").Value = -0.7907
$wsLosers.Range("E65").Value = 4.267

$wsLosers.Range("B66").Value = "ALLDIGI"
$wsLosers.Range("C66").Value = -2.4531
$wsLosers.Range("D66").Value = -0.045
$wsLosers.Range("E66").Value = -5.1342

$wsLosers.Range("B67").Value = "NSIL"
$wsLosers.Range("C67").Value = -2.4088
$wsLosers.Range("D67").Value = -1.7646
$wsLosers.Range("E67").Value = 4.7431

$wsLosers.Range("B68").Value = "COALINDIA"
$wsLosers.Range("C68").Value = -2.4016
$wsLosers.Range("D68").Value = -3.058
$wsLosers.Range("E68").Value = -2.0387

$wsLosers.Range("B69").Value = "FINOPB"
$wsLosers.Range("C69").Value = -2.3673
$wsLosers.Range("D69").Value = -6.2696
$wsLosers.Range("E69").Value = 11.1938

$wsLosers.Range("B70").Value = "UNIMECH"
$wsLosers.Range("C70").Value = -2.353
$wsLosers.Range("D70").Value = -1.1572
$wsLosers.Range("E70").Value = 0

$wsLosers.Range("B71").Value = "FCL"
$wsLosers.Range("C71").Value = -2.3453
$wsLosers.Range("D71").Value = -2.616
$wsLosers.Range("E71").Value = -0.02

$wsLosers.Range("B72").Value = "DEEDEV"
$wsLosers.Range("C72").Value = -2.3136
$wsLosers.Range("D72").Value = -6.6339
$wsLosers.Range("E72").Value = -7.4039

$wsLosers.Range("B73").Value = "WEALTH"
$wsLosers.Range("C73").Value = -2.3047
$wsLosers.Range("D73").Value = -3.8606
$wsLosers.Range("E73").Value = -2.8234

$wsLosers.Range("B74").Value = "RATNAMANI"
$wsLosers.Range("C74").Value = -2.2788
$wsLosers.Range("D74").Value = -0.4626
$wsLosers.Range("E74").Value = 0.8712

$wsLosers.Range("B75").Value = "CSBBANK"
$wsLosers.Range("C75").Value = -2.2695
$wsLosers.Range("D75").Value = 2.3137
$wsLosers.Range("E75").Value = 10.6999

$wsLosers.Range("B76").Value = "BBOX"
$wsLosers.Range("C76").Value = -2.2639
$wsLosers.Range("D76").Value = -4.7636
$wsLosers.Range("E76").Value = 5.1528

# --- 1 Month Performance sheet ---
$wsPerf = $wb.Worksheets.Item("1 Month Performance")

# Row 6: refreshed % change (IFBAGRO)
$wsPerf.Range("C6").Value = 65.2534

# Rows 36-40: refreshed list (values shift down one row, new entry IFBIND added at top)
$wsPerf.Range("B36").Value = "IFBIND"
$wsPerf.Range("C36").Value = 27.064

$wsPerf.Range("B37").Value = "MINDTECK"
$wsPerf.Range("C37").Value = 26.9415

$wsPerf.Range("B38").Value = "BHARATWIRE"
$wsPerf.Range("C38").Value = 26.5276

$wsPerf.Range("B39").Value = "HATSUN"
$wsPerf.Range("C39").Value = 26.492

$wsPerf.Range("B40").Value = "INDORAMA"
$wsPerf.Range("C40").Value = 26.4516
